$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("check_link") to make room for "orders_count"
$ws.Columns.Item(5).Insert()

# Set the new header
$ws.Cells.Item(1, 5).Value = "orders_count"

# Add the new data row (row 2) - force text format so numeric-looking
# strings (leading zeros, plain digits) are kept as text, not coerced to numbers
$rowRange = $ws.Range("A2:H2")
$rowRange.NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "0000001"
$ws.Cells.Item(2, 2).Value = "test"
$ws.Cells.Item(2, 3).Value = "test, test, test"
$ws.Cells.Item(2, 4).Value = "test"
$ws.Cells.Item(2, 5).Value = "3"
$ws.Cells.Item(2, 6).Value = "https://gift-bot-checks.s3.us-east-1.amazonaws.com/checks/AQADFOgxG3JVYEt-.jpg"
$ws.Cells.Item(2, 7).Value = "2024-12-25 22:28:44"
$ws.Cells.Item(2, 8).Value = "wakeupkstnv"
